$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the extra data rows (3-5) that were exported previously
$ws.Rows("3:5").Delete()

# Add a styled (underlined), empty placeholder cell at C12 and select it,
# mirroring the new export layout
$ws.Range("C12").Font.Underline = 1
$ws.Range("C12").Select()
